$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet stores Coin/Link/Price/Volume as plain text, including values
# that look numeric (e.g. "1.00", "0.602"). Excel's COM layer auto-detects
# numeric-looking strings and converts them to real numbers when assigned
# via Range.Value, which would change the cell type away from text. Force
# the cell to a text format first, then reset the style back to the default
# "Normal" so no stray number-format style lingers on the cell afterwards.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "65.326.34"
Set-TextValue $ws.Range("E2") "  +0.62%  "
Set-TextValue $ws.Range("D3") "3.369.21"
Set-TextValue $ws.Range("E3") "  +0.79%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.22%  "
Set-TextValue $ws.Range("D5") "182.70"
Set-TextValue $ws.Range("E5") "  -0.05%  "
Set-TextValue $ws.Range("D6") "537.92"
Set-TextValue $ws.Range("E6") "  +0.66%  "
Set-TextValue $ws.Range("D7") "0.602"
Set-TextValue $ws.Range("E7") "  -1.13%  "
Set-TextValue $ws.Range("D8") "3.361.75"
Set-TextValue $ws.Range("E8") "  +0.72%  "
Set-TextValue $ws.Range("E9") "  +0.10%  "
Set-TextValue $ws.Range("D10") "0.626"
Set-TextValue $ws.Range("E10") "  +1.33%  "
Set-TextValue $ws.Range("D11") "55.51"
Set-TextValue $ws.Range("E11") "  -7.52%  "
Set-TextValue $ws.Range("E12") "  +3.50%  "
Set-TextValue $ws.Range("D13") "0.0000265"
Set-TextValue $ws.Range("E13") "  +1.00%  "
Set-TextValue $ws.Range("D14") "9.21"
Set-TextValue $ws.Range("E14") "  -0.10%  "
Set-TextValue $ws.Range("D15") "3.921.75"
Set-TextValue $ws.Range("E15") "  +1.69%  "
Set-TextValue $ws.Range("B16") "WrappedEther"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D16") "3.382.93"
Set-TextValue $ws.Range("E16") "  +1.72%  "
Set-TextValue $ws.Range("B17") "TRON"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D17") "0.120"
Set-TextValue $ws.Range("E17") "  +1.24%  "
Set-TextValue $ws.Range("D18") "18.05"
Set-TextValue $ws.Range("E18") "  +1.70%  "
Set-TextValue $ws.Range("D19") "65.582.28"
Set-TextValue $ws.Range("E19") "  +1.27%  "
Set-TextValue $ws.Range("D20") "11.42"
Set-TextValue $ws.Range("E20") "  +0.91%  "
Set-TextValue $ws.Range("D21") "0.983"
Set-TextValue $ws.Range("E21") "  +0.94%  "
Set-TextValue $ws.Range("D22") "389.64"
Set-TextValue $ws.Range("E22") "  +2.60%  "
Set-TextValue $ws.Range("D23") "11.94"
Set-TextValue $ws.Range("E23") "  +4.37%  "
Set-TextValue $ws.Range("D24") "4.20"
Set-TextValue $ws.Range("E24") "  +5.88%  "
Set-TextValue $ws.Range("D25") "82.97"
Set-TextValue $ws.Range("E25") "  +1.89%  "
Set-TextValue $ws.Range("D26") "3.80"
Set-TextValue $ws.Range("E26") "  -1.33%  "
Set-TextValue $ws.Range("B27") "LEO"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D27") "6.13"
Set-TextValue $ws.Range("E27") "  +0.38%  "
Set-TextValue $ws.Range("B28") "ImmutableX"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D28") "2.82"
Set-TextValue $ws.Range("E28") "  +3.80%  "
Set-TextValue $ws.Range("D29") "11.52"
Set-TextValue $ws.Range("E29") "  -1.03%  "
Set-TextValue $ws.Range("D30") "8.45"
Set-TextValue $ws.Range("E30") "  -0.69%  "
Set-TextValue $ws.Range("D31") "29.47"
Set-TextValue $ws.Range("E31") "  +0.55%  "
Set-TextValue $ws.Range("D32") "659.59"
Set-TextValue $ws.Range("E32") "  -0.87%  "
Set-TextValue $ws.Range("D33") "6.77"
Set-TextValue $ws.Range("E33") "  -0.50%  "
Set-TextValue $ws.Range("D34") "11.43"
Set-TextValue $ws.Range("E34") "  +0.11%  "
Set-TextValue $ws.Range("E35") "  +0.58%  "
Set-TextValue $ws.Range("D36") "57.83"
Set-TextValue $ws.Range("E36") "  -3.48%  "
Set-TextValue $ws.Range("D37") "37.73"
Set-TextValue $ws.Range("E37") "  +1.22%  "
Set-TextValue $ws.Range("B38") "Dai"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D38") "1.00"
Set-TextValue $ws.Range("E38") "  +0.01%  "
Set-TextValue $ws.Range("B39") "TheGraph"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws.Range("D39") "0.398"
Set-TextValue $ws.Range("E39") "  -0.36%  "
Set-TextValue $ws.Range("D40") "0.0₃0772"
Set-TextValue $ws.Range("E40") "  +8.08%  "
Set-TextValue $ws.Range("D41") "2.79"
Set-TextValue $ws.Range("E41") "  +8.21%  "
Set-TextValue $ws.Range("B42") "FirstDigitalUSD"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D42") "1.00"
Set-TextValue $ws.Range("E42") "  +0.60%  "
Set-TextValue $ws.Range("B43") "Stacks"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D43") "3.27"
Set-TextValue $ws.Range("E43") "  +14.70%  "
Set-TextValue $ws.Range("D44") "0.129"
Set-TextValue $ws.Range("E44") "  +0.99%  "
Set-TextValue $ws.Range("D45") "3.024.64"
Set-TextValue $ws.Range("E45") "  +3.05%  "
Set-TextValue $ws.Range("E46") "  +0.93%  "
Set-TextValue $ws.Range("D47") "0.0411"
Set-TextValue $ws.Range("E47") "  +1.54%  "
Set-TextValue $ws.Range("D48") "2.71"
Set-TextValue $ws.Range("E48") "  +1.72%  "
Set-TextValue $ws.Range("D49") "3.19"
Set-TextValue $ws.Range("E49") "  +1.71%  "
Set-TextValue $ws.Range("D50") "0.127"
Set-TextValue $ws.Range("E50") "  -0.33%  "
Set-TextValue $ws.Range("D51") "8.53"
Set-TextValue $ws.Range("E51") "  +6.15%  "
